# Add "Direct CPI" and "Calculated CPI" columns (G, H) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same formatting as the other header cells (e.g. F1) to the new headers.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Header row text
$ws.Cells.Item(1, 7).Value = "Direct CPI"
$ws.Cells.Item(1, 8).Value = "Calculated CPI"

# Data rows: row number, Direct CPI, Calculated CPI
$data = @(
    @(2, 5.036522, 1.042728),
    @(3, 4.160849, 1.021555),
    @(4, 3.905553, 1.033842),
    @(5, 3.900594, 1.035275),
    @(6, 3.897679, 1.036175),
    @(7, 3.896822, 1.036477),
    @(8, 3.896752, 1.036492),
    @(9, 3.896616, 1.036555),
    @(10, 4.753193, 1.074854),
    @(11, 4.178524, 1.03553),
    @(12, 3.98374, 1.027977),
    @(13, 3.900594, 1.035275),
    @(14, 3.876813, 1.044055),
    @(15, 3.849235, 1.07351),
    @(16, 3.842645, 1.090709),
    @(17, 3.831575, 1.155634),
    @(18, 3.900918, 1.035878),
    @(19, 3.900594, 1.035275),
    @(20, 3.900594, 1.035275),
    @(21, 3.900594, 1.035275),
    @(22, 3.900594, 1.035275),
    @(23, 3.900594, 1.035275),
    @(24, 3.900594, 1.035275),
    @(25, 3.900594, 1.035275),
    @(26, 4.168985, 1.065469),
    @(27, 4.026669, 1.063076),
    @(28, 3.932666, 1.051251),
    @(29, 3.900594, 1.035275),
    @(30, 3.905948, 1.021406),
    @(31, 3.947798, 1.016099),
    @(32, 4.105177, 1.023663),
    @(33, 3.905159, 1.033921),
    @(34, 3.900594, 1.035275),
    @(35, 3.900079, 1.035429),
    @(36, 3.899845, 1.035474),
    @(37, 4.370536, 1.04594),
    @(38, 3.900594, 1.035275),
    @(39, 3.855532, 1.06282),
    @(40, 3.853747, 1.065393),
    @(41, 3.900594, 1.035275),
    @(42, 3.900594, 1.035275),
    @(43, 3.900594, 1.035275),
    @(44, 3.900594, 1.035275),
    @(45, 26.024219, 2.119923),
    @(46, 26.037528, 2.120258),
    @(47, 26.030926, 2.120666),
    @(48, 26.019584, 2.121013),
    @(49, 26.014169, 2.121198),
    @(50, 26.006216, 2.121268),
    @(51, 26.033599, 2.12128),
    @(52, 26.0345, 2.121282),
    @(53, 26.10692, 2.100689),
    @(54, 26.034201, 2.116132),
    @(55, 26.042438, 2.119731),
    @(56, 26.019584, 2.121013),
    @(57, 26.013784, 2.121162),
    @(58, 26.057559, 2.121165),
    @(59, 26.057559, 2.121166),
    @(60, 26.057559, 2.121168),
    @(61, 26.060408, 2.121013),
    @(62, 26.010508, 2.121013),
    @(63, 26.0671, 2.121013),
    @(64, 26.019584, 2.121013),
    @(65, 25.987687, 2.121012),
    @(66, 25.92413, 2.121012),
    @(67, 25.807272, 2.121009),
    @(68, 22.24489, 1.939559),
    @(69, 67.242693, 2.78396),
    @(70, 84.26831300000001, 2.485646),
    @(71, 45.530116, 2.242712),
    @(72, 26.019584, 2.121013),
    @(73, 16.903227, 2.059153),
    @(74, 12.882278, 2.023028),
    @(75, 9.523178, 1.934502),
    @(76, 26.010164, 2.120884),
    @(77, 26.019584, 2.121013),
    @(78, 26.029681, 2.121122),
    @(79, 26.030053, 2.121241),
    @(80, 26.040118, 2.118625),
    @(81, 26.019584, 2.121013),
    @(82, 26.016276, 2.121156),
    @(83, 26.016276, 2.121157),
    @(84, 25.994884, 2.121013),
    @(85, 26.019584, 2.121013),
    @(86, 26.058949, 2.121013),
    @(87, 26.027728, 2.121013)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $directCpi = $entry[1]
    $calcCpi = $entry[2]
    $ws.Cells.Item($row, 7).Value = $directCpi
    $ws.Cells.Item($row, 8).Value = $calcCpi
}
